# Auto-generated edit script: updates market price / profit columns (H-N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR leve-profit tables
# to reflect refreshed market board data from the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 619.625
$ws.Range("I38").Value = 55.833332
$ws.Range("J38").Value = 957.9
$ws.Range("K38").Value = 167.499996
$ws.Range("L38").Value = 2873.7
$ws.Range("M38").Value = 204.500004
$ws.Range("N38").Value = -3617.7

$ws.Range("H39").Value = 619.5
$ws.Range("I39").Value = 110
$ws.Range("J39").Value = 874.25
$ws.Range("K39").Value = 330
$ws.Range("L39").Value = 2622.75
$ws.Range("M39").Value = -34
$ws.Range("N39").Value = -3214.75

$ws.Range("H40").Value = 1469.6154
$ws.Range("I40").Value = 1603.3334
$ws.Range("J40").Value = 1429.5
$ws.Range("K40").Value = 1603.3334
$ws.Range("L40").Value = 1429.5
$ws.Range("M40").Value = -1428.3334
$ws.Range("N40").Value = -1779.5

$ws.Range("H52").Value = 120
$ws.Range("I52").Value = 120
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 360
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = -200
$ws.Range("N52").ClearContents()

$ws.Range("H74").Value = 74927.125
$ws.Range("I74").Value = 97414.37
$ws.Range("K74").Value = 97414.37
$ws.Range("M74").Value = -96478.37

$ws.Range("H77").Value = 74927.125
$ws.Range("I77").Value = 97414.37
$ws.Range("K77").Value = 487071.85
$ws.Range("M77").Value = -482391.85

$ws.Range("H116").Value = 12322.5
$ws.Range("J116").Value = 9784.286
$ws.Range("L116").Value = 9784.286
$ws.Range("N116").Value = -16668.286

$ws.Range("H137").Value = 1777.4062
$ws.Range("I137").Value = 1591.2084
$ws.Range("K137").Value = 4773.6252
$ws.Range("M137").Value = -2223.6252

$ws.Range("H138").Value = 2966.7424
$ws.Range("I138").Value = 1306.3125
$ws.Range("J138").Value = 4529.5
$ws.Range("K138").Value = 3918.9375
$ws.Range("L138").Value = 13588.5
$ws.Range("M138").Value = 1221.0625
$ws.Range("N138").Value = -23868.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4488.087
$ws.Range("I2").Value = 3621.1667
$ws.Range("J2").Value = 7609
$ws.Range("K2").Value = 3621.1667
$ws.Range("L2").Value = 7609
$ws.Range("M2").Value = -3508.1667
$ws.Range("N2").Value = -7835

$ws.Range("H61").Value = 3751.394
$ws.Range("I61").Value = 3348.8462
$ws.Range("J61").Value = 5246.5713
$ws.Range("K61").Value = 3348.8462
$ws.Range("L61").Value = 5246.5713
$ws.Range("M61").Value = -3136.8462
$ws.Range("N61").Value = -5670.5713

$ws.Range("H69").Value = 500299.66
$ws.Range("J69").Value = 500299.66
$ws.Range("L69").Value = 500299.66
$ws.Range("N69").Value = -501797.66

$ws.Range("H72").Value = 500299.66
$ws.Range("J72").Value = 500299.66
$ws.Range("L72").Value = 1500898.98
$ws.Range("N72").Value = -1508386.98

$ws.Range("H116").Value = 4488.087
$ws.Range("I116").Value = 3621.1667
$ws.Range("J116").Value = 7609
$ws.Range("K116").Value = 3621.1667
$ws.Range("L116").Value = 7609
$ws.Range("M116").Value = -1327.1667
$ws.Range("N116").Value = -12197

$ws.Range("H132").Value = 3629.6191
$ws.Range("I132").Value = 1782.9166
$ws.Range("K132").Value = 5348.7498
$ws.Range("M132").Value = -2818.7498

$ws.Range("H133").Value = 211499.5
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()

$ws.Range("H136").Value = 3751.394
$ws.Range("I136").Value = 3348.8462
$ws.Range("J136").Value = 5246.5713
$ws.Range("K136").Value = 10046.5386
$ws.Range("L136").Value = 15739.7139
$ws.Range("M136").Value = -7496.5386
$ws.Range("N136").Value = -20839.7139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4488.087
$ws.Range("I3").Value = 3621.1667
$ws.Range("J3").Value = 7609
$ws.Range("K3").Value = 3621.1667
$ws.Range("L3").Value = 7609
$ws.Range("M3").Value = -3507.1667
$ws.Range("N3").Value = -7837

$ws.Range("H62").Value = 80000
$ws.Range("J62").Value = 80000
$ws.Range("L62").Value = 80000
$ws.Range("N62").Value = -81372

$ws.Range("H65").Value = 80000
$ws.Range("J65").Value = 80000
$ws.Range("L65").Value = 240000
$ws.Range("N65").Value = -246864

$ws.Range("H86").Value = 3686.5625
$ws.Range("I86").Value = 2564.2222
$ws.Range("J86").Value = 5129.5713
$ws.Range("K86").Value = 2564.2222
$ws.Range("L86").Value = 5129.5713
$ws.Range("M86").Value = -1441.2222
$ws.Range("N86").Value = -7375.5713

$ws.Range("H89").Value = 3686.5625
$ws.Range("I89").Value = 2564.2222
$ws.Range("J89").Value = 5129.5713
$ws.Range("K89").Value = 12821.111
$ws.Range("L89").Value = 25647.8565
$ws.Range("M89").Value = -7205.111000000001
$ws.Range("N89").Value = -36879.85649999999

$ws.Range("H94").Value = 3287.5334
$ws.Range("I94").Value = 3342.8
$ws.Range("J94").Value = 3177
$ws.Range("K94").Value = 3342.8
$ws.Range("L94").Value = 3177
$ws.Range("M94").Value = -2891.8
$ws.Range("N94").Value = -4079

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 49257.684
$ws.Range("I31").Value = 56865.89
$ws.Range("K31").Value = 56865.89
$ws.Range("M31").Value = -56570.89

$ws.Range("H34").Value = 49257.684
$ws.Range("I34").Value = 56865.89
$ws.Range("K34").Value = 56865.89
$ws.Range("M34").Value = -56663.89

$ws.Range("H58").Value = 2374.2896
$ws.Range("I58").Value = 2452.2727
$ws.Range("K58").Value = 2452.2727
$ws.Range("M58").Value = -2249.2727

$ws.Range("H99").Value = 5558.88
$ws.Range("I99").Value = 4416.467
$ws.Range("J99").Value = 7272.5
$ws.Range("K99").Value = 4416.467
$ws.Range("L99").Value = 7272.5
$ws.Range("M99").Value = -2918.467
$ws.Range("N99").Value = -10268.5

$ws.Range("H126").Value = 5558.88
$ws.Range("I126").Value = 4416.467
$ws.Range("J126").Value = 7272.5
$ws.Range("K126").Value = 13249.401
$ws.Range("L126").Value = 21817.5
$ws.Range("M126").Value = -10779.401
$ws.Range("N126").Value = -26757.5

$ws.Range("H132").Value = 2934.15
$ws.Range("I132").Value = 3071.879
$ws.Range("K132").Value = 9215.636999999999
$ws.Range("M132").Value = -6685.636999999999

$ws.Range("H134").Value = 26081.285
$ws.Range("I134").Value = 20813.9
$ws.Range("K134").Value = 62441.7
$ws.Range("M134").Value = -59906.7

$ws.Range("H136").Value = 2374.2896
$ws.Range("I136").Value = 2452.2727
$ws.Range("K136").Value = 7356.8181
$ws.Range("M136").Value = -4806.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 138.18182
$ws.Range("J12").Value = 301.2
$ws.Range("L12").Value = 903.5999999999999
$ws.Range("N12").Value = -1249.6

$ws.Range("H68").Value = 6252837.5
$ws.Range("I68").Value = 2399.6667
$ws.Range("J68").Value = 10003100
$ws.Range("K68").Value = 7199.000100000001
$ws.Range("L68").Value = 30009300
$ws.Range("M68").Value = -6388.000100000001
$ws.Range("N68").Value = -30010922

$ws.Range("H71").Value = 6252837.5
$ws.Range("I71").Value = 2399.6667
$ws.Range("J71").Value = 10003100
$ws.Range("K71").Value = 21597.0003
$ws.Range("L71").Value = 90027900
$ws.Range("M71").Value = -17541.0003
$ws.Range("N71").Value = -90036012

$ws.Range("H132").Value = 1222.2
$ws.Range("I132").Value = 998.4286
$ws.Range("J132").Value = 1418
$ws.Range("K132").Value = 8985.857399999999
$ws.Range("L132").Value = 12762
$ws.Range("M132").Value = -6455.857399999999
$ws.Range("N132").Value = -17822

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3515.9473
$ws.Range("I132").Value = 3390
$ws.Range("K132").Value = 10170
$ws.Range("M132").Value = -7640

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2618.1936
$ws.Range("I132").Value = 2641.524
$ws.Range("J132").Value = 2569.2
$ws.Range("K132").Value = 7924.572
$ws.Range("L132").Value = 7707.599999999999
$ws.Range("M132").Value = -5394.572
$ws.Range("N132").Value = -12767.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 620.8421
$ws.Range("I107").Value = 661.6
$ws.Range("J107").Value = 575.55554
$ws.Range("K107").Value = 1984.8
$ws.Range("L107").Value = 1726.66662
$ws.Range("M107").Value = -64.80000000000018
$ws.Range("N107").Value = -5566.66662

$ws.Range("H126").Value = 2776.2778
$ws.Range("I126").Value = 2274.1428
$ws.Range("J126").Value = 4533.75
$ws.Range("K126").Value = 6822.428400000001
$ws.Range("L126").Value = 13601.25
$ws.Range("M126").Value = -4352.428400000001
$ws.Range("N126").Value = -18541.25

$ws.Range("H127").Value = 38466.332
$ws.Range("J127").Value = 38466.332
$ws.Range("L127").Value = 38466.332
$ws.Range("N127").Value = -48386.332

$ws.Range("H132").Value = 3487.7317
$ws.Range("I132").Value = 3377.2188
$ws.Range("K132").Value = 10131.6564
$ws.Range("M132").Value = -7601.6564

$ws.Range("H136").Value = 1598.4517
$ws.Range("I136").Value = 1584.48
$ws.Range("J136").Value = 1656.6666
$ws.Range("K136").Value = 4753.440000000001
$ws.Range("L136").Value = 4969.9998
$ws.Range("M136").Value = -2203.440000000001
$ws.Range("N136").Value = -10069.9998
